$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 3 de Septiembre de 2020 a las 12:10"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6291032
$ws.Range("C4").Value = 295
$ws.Range("D4").Value = 3547446
$ws.Range("E4").Value = 2553595
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 27
$ws.Range("H4").Value = 189991

# Row 17: Banglades
$ws.Range("A17").Value = "Banglades"
$ws.Range("B17").Value = 319686
$ws.Range("C17").Value = 2158
$ws.Range("D17").Value = 213980
$ws.Range("E17").Value = 101323
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 4383

# Row 29: Israel
$ws.Range("A29").Value = "Israel"
$ws.Range("B29").Value = 122779
$ws.Range("C29").Value = 1315
$ws.Range("D29").Value = 97885
$ws.Range("E29").Value = 23918
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 976

# Row 37: Rumania
$ws.Range("A37").Value = "Rumania"
$ws.Range("B37").Value = 91256
$ws.Range("C37").Value = 1365
$ws.Range("D37").Value = 39275
$ws.Range("E37").Value = 48216
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 44
$ws.Range("H37").Value = 3765

# Row 38: Kuwait
$ws.Range("A38").Value = "Kuwait"
$ws.Range("B38").Value = 87378
$ws.Range("C38").Value = 900
$ws.Range("D38").Value = 78791
$ws.Range("E38").Value = 8051
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 536

# Row 39: Oman
$ws.Range("A39").Value = "Oman"
$ws.Range("B39").Value = 86380
$ws.Range("C39").Value = 256
$ws.Range("D39").Value = 81828
$ws.Range("E39").Value = 3847
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 16
$ws.Range("H39").Value = 705

# Row 71: Austria
$ws.Range("A71").Value = "Austria"
$ws.Range("B71").Value = 28372
$ws.Range("C71").Value = 403
$ws.Range("D71").Value = 24171
$ws.Range("E71").Value = 3466
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 735

# Row 73: El Salvador
$ws.Range("A73").Value = "El Salvador"
$ws.Range("B73").Value = 26000
$ws.Range("C73").Value = 96
$ws.Range("D73").Value = 15119
$ws.Range("E73").Value = 10142
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 8
$ws.Range("H73").Value = 739

# Row 90: Croacia
$ws.Range("A90").Value = "Croacia"
$ws.Range("B90").Value = 11094
$ws.Range("C90").Value = 369
$ws.Range("D90").Value = 8266
$ws.Range("E90").Value = 2634
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 194

# Row 91: Noruega
$ws.Range("A91").Value = "Noruega"
$ws.Range("B91").Value = 11034
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 9348
$ws.Range("E91").Value = 1422
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 264

# Row 92: Grecia
$ws.Range("A92").Value = "Grecia"
$ws.Range("B92").Value = 10757
$ws.Range("C92").Value = 0
$ws.Range("D92").Value = 3804
$ws.Range("E92").Value = 6680
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 273

# Row 93: Consejo Danes para los Refugiados
$ws.Range("A93").Value = "Consejo Danes para los Refugiados"
$ws.Range("B93").Value = 10125
$ws.Range("C93").Value = 11
$ws.Range("D93").Value = 9367
$ws.Range("E93").Value = 499
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 259

# Row 96: Malasia
$ws.Range("A96").Value = "Malasia"
$ws.Range("B96").Value = 9374
$ws.Range("C96").Value = 14
$ws.Range("D96").Value = 9083
$ws.Range("E96").Value = 163
$ws.Range("F96").Value = 0
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 128

# Row 101: Finlandia
$ws.Range("A101").Value = "Finlandia"
$ws.Range("B101").Value = 8200
$ws.Range("C101").Value = 39
$ws.Range("D101").Value = 7350
$ws.Range("E101").Value = 514
$ws.Range("F101").Value = 0
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 336

# Row 112: Hong Kong
$ws.Range("A112").Value = "Hong Kong"
$ws.Range("B112").Value = 4839
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 4431
$ws.Range("E112").Value = 315
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 93

# Row 119: Eslovaquia
$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("B119").Value = 4163
$ws.Range("C119").Value = 121
$ws.Range("D119").Value = 2617
$ws.Range("E119").Value = 1509
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 4
$ws.Range("H119").Value = 37

# Row 120: Surinam
$ws.Range("A120").Value = "Surinam"
$ws.Range("B120").Value = 4149
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 3272
$ws.Range("E120").Value = 805
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 72

# Row 121: Cuba
$ws.Range("A121").Value = "Cuba"
$ws.Range("B121").Value = 4126
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 3458
$ws.Range("E121").Value = 570
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 98

# Row 122: Mozambique
$ws.Range("A122").Value = "Mozambique"
$ws.Range("B122").Value = 4117
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 2170
$ws.Range("E122").Value = 1922
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 25

# Row 123: Cabo Verde
$ws.Range("A123").Value = "Cabo Verde"
$ws.Range("B123").Value = 4048
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 3460
$ws.Range("E123").Value = 547
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 41

# Row 165: Vietnam
$ws.Range("A165").Value = "Vietnam"
$ws.Range("B165").Value = 1046
$ws.Range("C165").Value = 0
$ws.Range("D165").Value = 746
$ws.Range("E165").Value = 265
$ws.Range("F165").Value = 0
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 35

# Row 172: Islas Turcas y Caicos
$ws.Range("A172").Value = "Islas Turcas y Caicos"
$ws.Range("B172").Value = 555
$ws.Range("C172").Value = 17
$ws.Range("D172").Value = 220
$ws.Range("E172").Value = 331
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 4

# Row 179: Islas Feroe
$ws.Range("A179").Value = "Islas Feroe"
$ws.Range("B179").Value = 412
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 379
$ws.Range("E179").Value = 33
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

# Row 191: Brunei
$ws.Range("A191").Value = "Brunei"
$ws.Range("B191").Value = 145
$ws.Range("C191").Value = 1
$ws.Range("D191").Value = 139
$ws.Range("E191").Value = 3
$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 3
